$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed ticker cells for rows 2-11
$ws.Range("C2").Value = "NSE:AKZOINDIA"
$ws.Range("D2").Value = "NSE:ADANIENT"
$ws.Range("F2").Value = "NSE:ASTRAL"
$ws.Range("B3").Value = "NSE:ASTRAL"
$ws.Range("C3").Value = "NSE:BALMLAWRIE"
$ws.Range("D3").Value = "NSE:CUMMINSIND"
$ws.Range("E3").Value = "NSE:CONCOR"
$ws.Range("F3").Value = "NSE:GODREJPROP"
$ws.Range("B4").Value = "NSE:AUTOBEES"
$ws.Range("C4").Value = "NSE:GMDCLTD"
$ws.Range("D4").Value = "NSE:JKCEMENT"
$ws.Range("F4").Value = "NSE:INDHOTEL"
$ws.Range("B5").Value = "NSE:BANKETF"
$ws.Range("C5").Value = "NSE:GULPOLY"
$ws.Range("B6").Value = "NSE:GODREJPROP"
$ws.Range("C6").Value = "NSE:HCL-INSYS"
$ws.Range("B7").Value = "NSE:HDFCPVTBAN"
$ws.Range("C7").Value = "NSE:INTELLECT"
$ws.Range("B8").Value = "NSE:INDHOTEL"
$ws.Range("C8").Value = "NSE:JOCIL"
$ws.Range("B9").Value = "NSE:KALYANKJIL"
$ws.Range("C9").Value = "NSE:KRBL"
$ws.Range("B10").Value = "NSE:MANYAVAR"
$ws.Range("C10").Value = "NSE:PITTIENG"
$ws.Range("B11").Value = "NSE:NIFTYBEES"
$ws.Range("C11").Value = "NSE:SAKSOFT"

# Clear cells that become empty
$ws.Range("D5").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("F6").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("F9").ClearContents()
$ws.Range("F10").ClearContents()

# Remove now-unused rows 12-18 (table shrank from 16 to 9 data rows)
$ws.Range("A12:F18").EntireRow.Delete()

